# Generate Report for Archive
# Replace the "Ready for handoff" status text with "In Translation"
# across the Overview sheet (zh-cn / de-de columns) and the per-language
# status sheets (zh-cn, de-de), matching the author's commit.

$wb = $excel.ActiveWorkbook

$oldText = "Ready for handoff"
$newText = "In Translation"

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count
    $colCount = $usedRange.Columns.Count

    for ($r = 1; $r -le $rowCount; $r++) {
        for ($c = 1; $c -le $colCount; $c++) {
            $cell = $usedRange.Cells.Item($r, $c)
            # NOTE: compare with the literal on the left - $cell.Value2 can come
            # back as a native boolean for True/False cells, and PowerShell's
            # "-eq" coerces the right-hand side to the left operand's type,
            # which would make "$true -eq 'any non-empty string'" true.
            if ($oldText -eq $cell.Value2) {
                $cell.Value = $newText
            }
        }
    }
}

# The shortened status text means the Status-related columns (zh-cn / de-de on
# the Overview sheet, and the Status column on each language sheet) now need
# less horizontal space, so narrow them accordingly.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth   # zh-cn status column
$overview.Columns.Item(6).ColumnWidth = $newWidth   # de-de status column

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth        # Status column

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth        # Status column
